# Generate Report for Handoff
# Updates the localization-status workbook to reflect that the handoff
# package is now ready (status text + new handoff generation timestamps),
# and widens the Status columns to fit the new, longer text.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

# --- Overview sheet ---
$overview.Range("E2").Value = "Ready for handoff"
$overview.Range("F2").Value = "Ready for handoff"
$overview.Range("G2").Value = "2016-08-22 16:41:39"

# --- zh-cn sheet ---
$zhcn.Range("C2").Value = "Ready for handoff"
$zhcn.Range("H2").Value = "2016-08-22 16:41:34"

# --- de-de sheet ---
$dede.Range("C2").Value = "Ready for handoff"
$dede.Range("H2").Value = "2016-08-22 16:41:39"

# Widen the Status columns so the new "Ready for handoff" text fits.
# (ColumnWidth is quantized to Excel's pixel grid on save, so we feed it
# the character-width value that rounds closest to the target 17.22.)
$newStatusWidth = 98 / 6
$overview.Range("E1").ColumnWidth = $newStatusWidth
$overview.Range("F1").ColumnWidth = $newStatusWidth
$zhcn.Range("C1").ColumnWidth = $newStatusWidth
$dede.Range("C1").ColumnWidth = $newStatusWidth
